$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  C="110111111111101110010101100111"; D=0.7655070970350591;  E=0.0004302871702627556; F=0.2192471784263435 },
    @{ Row=3;  C="110111111111110101111010110010"; D=0.7655577013447894;  E=0.04022746809695271;   F=0.5377540570379942 },
    @{ Row=4;  C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.5768111863966354;    F=0.7277932274751817 },
    @{ Row=5;  C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.5768111863966354;    F=0.7466779934933582 },
    @{ Row=6;  C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577013447894;    F=0.765557707863761  },
    @{ Row=7;  C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=8;  C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=9;  C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=10; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=11; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=12; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=13; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=14; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=15; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=16; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=17; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=18; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=19; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=20; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 },
    @{ Row=21; C="110111111111110101111010110111"; D=0.7655577094935038;  E=0.7655577094935038;    F=0.7655577094935037 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row

    # Column C holds a chromosome bit-string that must stay text (it looks
    # like a huge number and would otherwise be auto-coerced). Force text
    # formatting just long enough to assign it, then restore the default
    # "Normal" style so no stray number-format style sticks to the cell.
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $rowData.C
    $cCell.Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
}
